$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (F column) for three events
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F3").Value = 351   # 南宁·布谷鸟动漫展5th: 348 -> 351
$wsExhibit.Range("F4").Value = 4706  # 南宁·2024良牙动漫秋季盛典（秋典）: 4680 -> 4706
$wsExhibit.Range("F6").Value = 475   # 南宁·万圣漫控嘉年华10: 473 -> 475

# Sheet "全部类型" (sheet4): same three events, different row for the third one
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 351   # 南宁·布谷鸟动漫展5th: 348 -> 351
$wsAll.Range("F4").Value = 4706  # 南宁·2024良牙动漫秋季盛典（秋典）: 4680 -> 4706
$wsAll.Range("F8").Value = 475   # 南宁·万圣漫控嘉年华10: 473 -> 475
